# Update the Directors roster for Summer21 ("Updated Ec and Directors"):
# - The officer list shrinks from 12 rows to 10 rows (drop the trailing two rows)
# - Swap in the new slate of names, keeping "Adora Chen" (row 5) in place,
#   and reassign titles for every director row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer needed (old rows 11 & 12)
$ws.Range("A11:C12").EntireRow.Delete()

# Replace names (column A) - "Adora Chen" on row 5 is unchanged, so it is skipped.
# Order below matches how the names were actually retyped.
$ws.Cells.Item(4, 1).Value  = "Kevin Cao"
$ws.Cells.Item(2, 1).Value  = "Jacob Won"
$ws.Cells.Item(3, 1).Value  = "Mirsab Mirza"
$ws.Cells.Item(6, 1).Value  = "Kelsie Kim"
$ws.Cells.Item(7, 1).Value  = "Dylan Tanzil"
$ws.Cells.Item(8, 1).Value  = "Yoyo Cao"
$ws.Cells.Item(9, 1).Value  = "David Ayala"
$ws.Cells.Item(10, 1).Value = "Safah Faraz"

# Update titles (column B) to match each director's new role
$ws.Cells.Item(2, 2).Value  = "Marketing"
$ws.Cells.Item(3, 2).Value  = "Professional Development"
$ws.Cells.Item(4, 2).Value  = "Technology"
$ws.Cells.Item(5, 2).Value  = "Content Creation"
$ws.Cells.Item(6, 2).Value  = "Brotherhood and Social Activities"
$ws.Cells.Item(7, 2).Value  = "Brotherhood and Social Activities"
$ws.Cells.Item(8, 2).Value  = "University Relations"
$ws.Cells.Item(9, 2).Value  = "External Relations"
$ws.Cells.Item(10, 2).Value = "Fundraising"

$ws.Range("B11").Select()
